$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44372
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 55
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "`$/caja 15 kilos granel"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 15

$ws.Range("D3").Value = 44372
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 70
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("Q3").Value = "`$/caja 15 kilos granel"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 800
$ws.Range("T3").Value = 15

$ws.Range("D4").Value = 44327
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 4
$ws.Range("N4").Value = 150000
$ws.Range("O4").Value = 150000
$ws.Range("P4").Value = 150000
$ws.Range("Q4").Value = "`$/bins (450 kilos)"
$ws.Range("R4").Value = "Provincia de Cachapoal"
$ws.Range("S4").Value = 333
$ws.Range("T4").Value = 450

$ws.Range("D5").Value = 44364
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = "`$/caja 15 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 800
$ws.Range("T5").Value = 15

$ws.Range("D6").Value = 44364
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("Q6").Value = "`$/caja 15 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 667
$ws.Range("T6").Value = 15

$ws.Range("D7").Value = 44364
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("Q7").Value = "`$/caja 15 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 533
$ws.Range("T7").Value = 15

$ws.Range("D8").Value = 44376
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = "`$/caja 15 kilos granel"
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 15

$ws.Range("D9").Value = 44376
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 85
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("Q9").Value = "`$/caja 15 kilos granel"
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 800
$ws.Range("T9").Value = 15

$ws.Range("D10").Value = 44307
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 6
$ws.Range("N10").Value = 250000
$ws.Range("O10").Value = 250000
$ws.Range("P10").Value = 250000
$ws.Range("Q10").Value = "`$/bins (450 kilos)"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 556
$ws.Range("T10").Value = 450

$ws.Range("D11").Value = 44405
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 16000
$ws.Range("Q11").Value = "`$/caja 18 kilos granel"
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 889
$ws.Range("T11").Value = 18

$ws.Range("D12").Value = 44405
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = 12500
$ws.Range("O12").Value = 12500
$ws.Range("P12").Value = 12500
$ws.Range("Q12").Value = "`$/caja 18 kilos granel"
$ws.Range("R12").Value = "Región Metropolitana"
$ws.Range("S12").Value = 694
$ws.Range("T12").Value = 18

$ws.Range("D13").Value = 44384
$ws.Range("L13").Value = "Especial"
$ws.Range("M13").Value = 70
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 14000
$ws.Range("Q13").Value = "`$/caja 15 kilos granel"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 933
$ws.Range("T13").Value = 15

$ws.Range("D14").Value = 44384
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("Q14").Value = "`$/caja 15 kilos granel"
$ws.Range("R14").Value = "Región de O'Higgins"
$ws.Range("S14").Value = 800
$ws.Range("T14").Value = 15

$ws.Range("D15").Value = 44384
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("Q15").Value = "`$/caja 15 kilos granel"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 667
$ws.Range("T15").Value = 15

$ws.Range("D16").Value = 44383
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 70
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 14000
$ws.Range("Q16").Value = "`$/caja 15 kilos granel"
$ws.Range("R16").Value = "Región de O'Higgins"
$ws.Range("S16").Value = 933
$ws.Range("T16").Value = 15

$ws.Range("D17").Value = 44383
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("Q17").Value = "`$/caja 15 kilos granel"
$ws.Range("R17").Value = "Región de O'Higgins"
$ws.Range("S17").Value = 800
$ws.Range("T17").Value = 15

$ws.Range("D18").Value = 44292
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 10500
$ws.Range("O18").Value = 11000
$ws.Range("P18").Value = 10775
$ws.Range("Q18").Value = "`$/caja 18 kilos granel"
$ws.Range("R18").Value = "Región de O'Higgins"
$ws.Range("S18").Value = 599
$ws.Range("T18").Value = 18

$ws.Range("D19").Value = 44314
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 7
$ws.Range("N19").Value = 230000
$ws.Range("O19").Value = 230000
$ws.Range("P19").Value = 230000
$ws.Range("Q19").Value = "`$/bins (450 kilos)"
$ws.Range("R19").Value = "Región Metropolitana"
$ws.Range("S19").Value = 511
$ws.Range("T19").Value = 450

$ws.Range("D20").Value = 44363
$ws.Range("L20").Value = "Especial"
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 12000
$ws.Range("O20").Value = 12000
$ws.Range("P20").Value = 12000
$ws.Range("Q20").Value = "`$/caja 15 kilos granel"
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 800
$ws.Range("T20").Value = 15

$ws.Range("D21").Value = 44363
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 120
$ws.Range("N21").Value = 10000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 10000
$ws.Range("Q21").Value = "`$/caja 15 kilos granel"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 667
$ws.Range("T21").Value = 15

$ws.Range("D22").Value = 44363
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 80
$ws.Range("N22").Value = 8000
$ws.Range("O22").Value = 8000
$ws.Range("P22").Value = 8000
$ws.Range("Q22").Value = "`$/caja 15 kilos granel"
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 533
$ws.Range("T22").Value = 15

$ws.Range("D23").Value = 44316
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 4
$ws.Range("N23").Value = 230000
$ws.Range("O23").Value = 230000
$ws.Range("P23").Value = 230000
$ws.Range("Q23").Value = "`$/bins (450 kilos)"
$ws.Range("R23").Value = "Región Metropolitana"
$ws.Range("S23").Value = 511
$ws.Range("T23").Value = 450

$ws.Range("D24").Value = 44386
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 45
$ws.Range("N24").Value = 14000
$ws.Range("O24").Value = 14000
$ws.Range("P24").Value = 14000
$ws.Range("Q24").Value = "`$/caja 15 kilos granel"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 933
$ws.Range("T24").Value = 15

$ws.Range("D25").Value = 44386
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 60
$ws.Range("N25").Value = 12000
$ws.Range("O25").Value = 12000
$ws.Range("P25").Value = 12000
$ws.Range("Q25").Value = "`$/caja 15 kilos granel"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 800
$ws.Range("T25").Value = 15

$ws.Range("D26").Value = 44313
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 6
$ws.Range("N26").Value = 240000
$ws.Range("O26").Value = 240000
$ws.Range("P26").Value = 240000
$ws.Range("Q26").Value = "`$/bins (450 kilos)"
$ws.Range("R26").Value = "Región de O'Higgins"
$ws.Range("S26").Value = 533
$ws.Range("T26").Value = 450

$ws.Range("D27").Value = 44315
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 5
$ws.Range("N27").Value = 230000
$ws.Range("O27").Value = 230000
$ws.Range("P27").Value = 230000
$ws.Range("Q27").Value = "`$/bins (450 kilos)"
$ws.Range("R27").Value = "Región Metropolitana"
$ws.Range("S27").Value = 511
$ws.Range("T27").Value = 450

$ws.Range("D28").Value = 44355
$ws.Range("L28").Value = "Especial"
$ws.Range("M28").Value = 120
$ws.Range("N28").Value = 12000
$ws.Range("O28").Value = 12000
$ws.Range("P28").Value = 12000
$ws.Range("Q28").Value = "`$/caja 15 kilos granel"
$ws.Range("R28").Value = "Región Metropolitana"
$ws.Range("S28").Value = 800
$ws.Range("T28").Value = 15

$ws.Range("D29").Value = 44355
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 150
$ws.Range("N29").Value = 10000
$ws.Range("O29").Value = 10000
$ws.Range("P29").Value = 10000
$ws.Range("Q29").Value = "`$/caja 15 kilos granel"
$ws.Range("R29").Value = "Región Metropolitana"
$ws.Range("S29").Value = 667
$ws.Range("T29").Value = 15

$ws.Range("D30").Value = 44301
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 150
$ws.Range("N30").Value = 12000
$ws.Range("O30").Value = 12000
$ws.Range("P30").Value = 12000
$ws.Range("Q30").Value = "`$/caja 18 kilos granel"
$ws.Range("R30").Value = "Provincia de Cachapoal"
$ws.Range("S30").Value = 667
$ws.Range("T30").Value = 18

$ws.Range("D31").Value = 44301
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 10000
$ws.Range("O31").Value = 10000
$ws.Range("P31").Value = 10000
$ws.Range("Q31").Value = "`$/caja 18 kilos granel"
$ws.Range("R31").Value = "Provincia de Cachapoal"
$ws.Range("S31").Value = 556
$ws.Range("T31").Value = 18
